$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("Main ISA")
$wsR = $wb.Worksheets.Item("R-Type")

# Fill in the newly-designed/tested AluControl codes for NOR, XOR (rows 26-27)
# and the Shift operations (rows 2-7) on the R-Type sheet. Shared strings are
# created by the engine in the order the cells are written, so we write the
# NOR/XOR rows first to match the upstream commit's string order.
$wsR.Range("F26").Value = "4'b0011"
$wsR.Range("F27").Value = "4'b0100"

$wsR.Range("F2").Value = "4'b1000"
$wsR.Range("F3").Value = "4'b1001"
$wsR.Range("F4").Value = "4'b1010"
$wsR.Range("F5").Value = "4'b1011"
$wsR.Range("F6").Value = "4'b1100"
$wsR.Range("F7").Value = "4'b1101"

# Recreate the author's final UI state: Main ISA scrolled/selected at D10,
# then R-Type made the active (selected) sheet with F7 selected.
$wsMain.Activate()
$wsMain.Range("D10").Select()

$wsR.Activate()
$wsR.Range("F7").Select()
